# Generate Report for Handoff
# Updates the "d5289b5b-7386-4bf1-ae40-3302f0f8189c" row across the
# Overview / zh-cn / de-de sheets to reflect that the file is ready
# for handoff again (handback version is stale), and records the
# relevant timestamps + error detail message.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/934d28e12adf7f2230fe4afbe68b5bbbff36a7a4/e2e/d5289b5b-7386-4bf1-ae40-3302f0f8189c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e22b3384db98a183da3729b4f1e673f993989e0/e2e/d5289b5b-7386-4bf1-ae40-3302f0f8189c.md."

# --- Overview sheet: row 3 is the d5289b5b file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReadyForHandoff
$wsOverview.Range("F3").Value = $statusReadyForHandoff
$wsOverview.Range("G3").Value = "2016-09-06 05:01:51"

# --- zh-cn sheet: row 3 is the d5289b5b file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReadyForHandoff
$wsZhCn.Range("H3").Value = "2016-09-06 05:01:46"
$wsZhCn.Range("P3").Value = $errorDetail
# Widen the Error Detail column to fit the new message (matches width=40)
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the d5289b5b file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReadyForHandoff
$wsDeDe.Range("H3").Value = "2016-09-06 05:01:51"
$wsDeDe.Range("P3").Value = $errorDetail
# Widen the Error Detail column to fit the new message (matches width=40)
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
